# Auto-generated edit script: update crypto price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.646.28"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "1.632.52"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'213.15"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'0.494"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").Value = "'18.98"
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("D11").Value = "'0.0842"
$ws.Range("E11").Value = "  +3.51%  "
$ws.Range("D12").Value = "1.860.41"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "1.638.55"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "'0.526"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "26.647.64"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "'63.10"
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "'210.11"
$ws.Range("D21").Value = "'4.30"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  +3.00%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").Value = "'147.24"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").Value = "'6.91"
$ws.Range("E28").Value = "  +5.24%  "
$ws.Range("D29").Value = "'15.40"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "'0.0521"
$ws.Range("E30").Value = "  +4.69%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "'2.94"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "1.169.52"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'0.505"
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "1.771.09"
$ws.Range("E44").Value = "  +1.46%  "
$ws.Range("D45").Value = "'92.64"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "'54.56"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "'0.0513"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("E49").Value = "  +4.58%  "
$ws.Range("E51").Value = "  -0.05%  "
